$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counts and percentages for "Target Type" categorical summary
$ws.Range("B2").Value = 264
$ws.Range("C2").Value = 91.35

$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 8.65
